$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("player_parameter")

# Fix parameters
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 40
$ws.Range("B4").Value = 20
$ws.Range("B5").Value = 0.3
$ws.Range("B20").Value = 150

# Update selection to B3
$ws.Range("B3").Select()
